$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.996.72"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "'3.502.80"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'605.57"

$ws.Range("D6").Value = "'172.49"
$ws.Range("E6").Value = "  -0.49%  "

$ws.Range("E7").Value = "  -1.58%  "

$ws.Range("D8").Value = "'3.495.62"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "'0.194"
$ws.Range("E10").Value = "  -3.09%  "

$ws.Range("D11").Value = "'7.24"
$ws.Range("E11").Value = "  +7.55%  "

$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").Value = "'45.97"
$ws.Range("E13").Value = "  -2.99%  "

$ws.Range("E14").Value = "  -1.83%  "

$ws.Range("D15").Value = "'4.069.12"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "'8.34"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("D17").Value = "'613.32"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").Value = "'3.494.78"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("D19").Value = "'70.020.07"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").Value = "'17.49"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").Value = "'0.877"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").Value = "'9.12"
$ws.Range("E23").Value = "  -8.73%  "

$ws.Range("D24").Value = "'98.54"
$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").Value = "'15.51"
$ws.Range("E25").Value = "  -2.51%  "

$ws.Range("D26").Value = "'3.72"
$ws.Range("E26").Value = "  -3.56%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").Value = "'2.56"
$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("D30").Value = "'8.99"
$ws.Range("E30").Value = "  -2.61%  "

$ws.Range("D31").Value = "'2.98"
$ws.Range("E31").Value = "  -3.78%  "

$ws.Range("D32").Value = "'8.04"
$ws.Range("E32").Value = "  -5.29%  "

$ws.Range("E33").Value = "  -4.62%  "

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'633.02"
$ws.Range("E34").Value = "  +11.14%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.80"
$ws.Range("E35").Value = "  -2.92%  "

$ws.Range("D36").Value = "'0.0996"
$ws.Range("E36").Value = "  -2.34%  "

$ws.Range("D37").Value = "'10.73"
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("E38").Value = "  +4.55%  "

$ws.Range("D39").Value = "'3.46"
$ws.Range("E39").Value = "  -4.40%  "

$ws.Range("D40").Value = "'56.69"
$ws.Range("E40").Value = "  -1.77%  "

$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("E42").Value = "  +1.25%  "

$ws.Range("D43").Value = "'3.351.53"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "'0.0₃0732"
$ws.Range("E44").Value = "  +2.57%  "

$ws.Range("D45").Value = "'0.310"
$ws.Range("E45").Value = "  -5.69%  "

$ws.Range("D46").Value = "'2.90"
$ws.Range("E46").Value = "  -4.43%  "

$ws.Range("D47").Value = "'31.80"
$ws.Range("E47").Value = "  -3.93%  "

$ws.Range("D48").Value = "'2.55"
$ws.Range("E48").Value = "  -4.20%  "

$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").Value = "'133.29"
$ws.Range("E50").Value = "  -0.32%  "

$ws.Range("E51").Value = "  -0.03%  "
